# Weekly refresh of the Caqui (persimmon) price sheet:
# a new week's record is inserted at the top of the data block (row 9),
# pushing all the existing weekly records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 - shifts rows 9:39 down to 10:40
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with this week's record
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44676
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107001
$ws.Range("J9").Value = "Caqui"
$ws.Range("K9").Value = "Fuyu"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 115
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = "`$/bandeja 15 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 15
